$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.039913953218831
$ws.Cells.Item(2, 4).Value = 1.043047995616468
$ws.Cells.Item(2, 5).Value = 1.047803087651592
$ws.Cells.Item(2, 6).Value = 1.057294117141863
$ws.Cells.Item(2, 9).Value = 1.040920471246013
$ws.Cells.Item(2, 10).Value = 1.045003393771235
$ws.Cells.Item(2, 11).Value = 1.045822961556264
$ws.Cells.Item(2, 12).Value = 1.050564705429315
$ws.Cells.Item(2, 13).Value = 1.060029481313317
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.040778764767479
$ws.Cells.Item(3, 4).Value = 1.043702153493407
$ws.Cells.Item(3, 5).Value = 1.048583400728539
$ws.Cells.Item(3, 6).Value = 1.058181755642709
$ws.Cells.Item(3, 9).Value = 1.041126990835344
$ws.Cells.Item(3, 10).Value = 1.045514030462324
$ws.Cells.Item(3, 11).Value = 1.046288485527668
$ws.Cells.Item(3, 12).Value = 1.051157018058541
$ws.Cells.Item(3, 13).Value = 1.060730740038175
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.041338989996178
$ws.Cells.Item(4, 4).Value = 1.044125928433513
$ws.Cells.Item(4, 5).Value = 1.049089277680735
$ws.Cells.Item(4, 6).Value = 1.058757269054844
$ws.Cells.Item(4, 9).Value = 1.0412597713379
$ws.Cells.Item(4, 10).Value = 1.04584441981079
$ws.Cells.Item(4, 11).Value = 1.046589515701228
$ws.Cells.Item(4, 12).Value = 1.051540585455486
$ws.Cells.Item(4, 13).Value = 1.061185013231134
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.041574659011278
$ws.Cells.Item(5, 4).Value = 1.044304199282138
$ws.Cells.Item(5, 5).Value = 1.049302176952785
$ws.Cells.Item(5, 6).Value = 1.058999488836886
$ws.Cells.Item(5, 9).Value = 1.04131538772072
$ws.Cells.Item(5, 10).Value = 1.045983307862026
$ws.Cells.Item(5, 11).Value = 1.04671602087977
$ws.Cells.Item(5, 12).Value = 1.05170190798808
$ws.Cells.Item(5, 13).Value = 1.061376110671841
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.04161423764686
$ws.Cells.Item(6, 4).Value = 1.044334138487063
$ws.Cells.Item(6, 5).Value = 1.049337937031652
$ws.Cells.Item(6, 6).Value = 1.059040174614978
$ws.Cells.Item(6, 9).Value = 1.041324713941632
$ws.Cells.Item(6, 10).Value = 1.046006627293746
$ws.Cells.Item(6, 11).Value = 1.046737258810853
$ws.Cells.Item(6, 12).Value = 1.05172899886258
$ws.Cells.Item(6, 13).Value = 1.061408203822587
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.041342138428311
$ws.Cells.Item(7, 4).Value = 1.044128310044183
$ws.Cells.Item(7, 5).Value = 1.049092121555183
$ws.Cells.Item(7, 6).Value = 1.058760504533315
$ws.Cells.Item(7, 9).Value = 1.041260515291297
$ws.Cells.Item(7, 10).Value = 1.045846275671253
$ws.Cells.Item(7, 11).Value = 1.046591206259398
$ws.Cells.Item(7, 12).Value = 1.051542740778608
$ws.Cells.Item(7, 13).Value = 1.061187566211036
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.040206088067412
$ws.Cells.Item(8, 4).Value = 1.043268968639539
$ws.Cells.Item(8, 5).Value = 1.048066597821805
$ws.Cells.Item(8, 6).Value = 1.057593859107093
$ws.Cells.Item(8, 9).Value = 1.040990441395056
$ws.Cells.Item(8, 10).Value = 1.045175970627664
$ws.Cells.Item(8, 11).Value = 1.045980327190048
$ws.Cells.Item(8, 12).Value = 1.05076481699089
$ws.Cells.Item(8, 13).Value = 1.060266368050175
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.038209148168881
$ws.Cells.Item(9, 4).Value = 1.04175853508243
$ws.Cells.Item(9, 5).Value = 1.046266940983305
$ws.Cells.Item(9, 6).Value = 1.055546984556841
$ws.Cells.Item(9, 9).Value = 1.040508046116371
$ws.Cells.Item(9, 10).Value = 1.043994657162588
$ws.Cells.Item(9, 11).Value = 1.04490244150318
$ws.Cells.Item(9, 12).Value = 1.049396386447455
$ws.Cells.Item(9, 13).Value = 1.058647088860734
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.0368812579333
$ws.Cells.Item(10, 4).Value = 1.04075425943529
$ws.Cells.Item(10, 5).Value = 1.045072281360775
$ws.Cells.Item(10, 6).Value = 1.05418849736583
$ws.Cells.Item(10, 9).Value = 1.040182126833609
$ws.Cells.Item(10, 10).Value = 1.043207087088128
$ws.Cells.Item(10, 11).Value = 1.044182957158901
$ws.Cells.Item(10, 12).Value = 1.048485775847303
$ws.Cells.Item(10, 13).Value = 1.057570344927497
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.036307093053398
$ws.Cells.Item(11, 4).Value = 1.04032005401451
$ws.Cells.Item(11, 5).Value = 1.044556215049024
$ws.Cells.Item(11, 6).Value = 1.053601726352139
$ws.Cells.Item(11, 9).Value = 1.04003998368607
$ws.Cells.Item(11, 10).Value = 1.042866069399387
$ws.Cells.Item(11, 11).Value = 1.043871215947937
$ws.Cells.Item(11, 12).Value = 1.048091886438725
$ws.Cells.Item(11, 13).Value = 1.057104780858035
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.036093947425671
$ws.Cells.Item(12, 4).Value = 1.040158870554561
$ws.Cells.Item(12, 5).Value = 1.044364711502949
$ws.Cells.Item(12, 6).Value = 1.053383995123005
$ws.Cells.Item(12, 9).Value = 1.039987033144472
$ws.Cells.Item(12, 10).Value = 1.042739402285815
$ws.Cells.Item(12, 11).Value = 1.043755392612211
$ws.Cells.Item(12, 12).Value = 1.047945641521638
$ws.Cells.Item(12, 13).Value = 1.056931952348123
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.036139662215581
$ws.Cells.Item(13, 4).Value = 1.040193440411531
$ws.Cells.Item(13, 5).Value = 1.044405781197239
$ws.Cells.Item(13, 6).Value = 1.05343068915059
$ws.Cells.Item(13, 9).Value = 1.039998398093858
$ws.Cells.Item(13, 10).Value = 1.042766572699494
$ws.Cells.Item(13, 11).Value = 1.04378023839288
$ws.Cells.Item(13, 12).Value = 1.047977008672756
$ws.Cells.Item(13, 13).Value = 1.056969019976737
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.036289471815148
$ws.Cells.Item(14, 4).Value = 1.040306728494247
$ws.Cells.Item(14, 5).Value = 1.044540381492097
$ws.Cells.Item(14, 6).Value = 1.053583724088398
$ws.Cells.Item(14, 9).Value = 1.040035609879606
$ws.Cells.Item(14, 10).Value = 1.0428555990075
$ws.Cells.Item(14, 11).Value = 1.043861642534597
$ws.Cells.Item(14, 12).Value = 1.048079796491971
$ws.Cells.Item(14, 13).Value = 1.057090492694998
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.036381791033139
$ws.Cells.Item(15, 4).Value = 1.040376542292608
$ws.Cells.Item(15, 5).Value = 1.044623337937373
$ws.Cells.Item(15, 6).Value = 1.053678043391179
$ws.Cells.Item(15, 9).Value = 1.040058517135258
$ws.Cells.Item(15, 10).Value = 1.042910451361392
$ws.Cells.Item(15, 11).Value = 1.043911794537806
$ws.Cells.Item(15, 12).Value = 1.048143135869338
$ws.Cells.Item(15, 13).Value = 1.057165349703815
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.036919380779357
$ws.Cells.Item(16, 4).Value = 1.040783090115158
$ws.Cells.Item(16, 5).Value = 1.0451065570325
$ws.Cells.Item(16, 6).Value = 1.05422747040653
$ws.Cells.Item(16, 9).Value = 1.040191539017816
$ws.Cells.Item(16, 10).Value = 1.04322971952722
$ws.Cells.Item(16, 11).Value = 1.044203642293808
$ws.Cells.Item(16, 12).Value = 1.048511925768732
$ws.Cells.Item(16, 13).Value = 1.057601257199028
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.037256817331171
$ws.Cells.Item(17, 4).Value = 1.041038282751759
$ws.Cells.Item(17, 5).Value = 1.045409997955483
$ws.Cells.Item(17, 6).Value = 1.054572504435851
$ws.Cells.Item(17, 9).Value = 1.040274708081472
$ws.Cells.Item(17, 10).Value = 1.043429990317349
$ws.Cells.Item(17, 11).Value = 1.044386658078625
$ws.Cells.Item(17, 12).Value = 1.048743368951554
$ws.Cells.Item(17, 13).Value = 1.057874871793119
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.037453717412751
$ws.Cells.Item(18, 4).Value = 1.041187195146208
$ws.Cells.Item(18, 5).Value = 1.045587108390975
$ws.Cells.Item(18, 6).Value = 1.054773897971482
$ws.Cells.Item(18, 9).Value = 1.040323120925809
$ws.Cells.Item(18, 10).Value = 1.043546805314227
$ws.Cells.Item(18, 11).Value = 1.044493388734164
$ws.Cells.Item(18, 12).Value = 1.048878405405032
$ws.Cells.Item(18, 13).Value = 1.058034531400699
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.037520868614137
$ws.Cells.Item(19, 4).Value = 1.041237981031882
$ws.Cells.Item(19, 5).Value = 1.045647518514209
$ws.Cells.Item(19, 6).Value = 1.054842591825985
$ws.Cells.Item(19, 9).Value = 1.040339611766574
$ws.Cells.Item(19, 10).Value = 1.043586636245305
$ws.Cells.Item(19, 11).Value = 1.044529777808839
$ws.Cells.Item(19, 12).Value = 1.048924456013589
$ws.Cells.Item(19, 13).Value = 1.05808898215153
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.037220605395126
$ws.Cells.Item(20, 4).Value = 1.04101089648395
$ws.Cells.Item(20, 5).Value = 1.045377429348448
$ws.Cells.Item(20, 6).Value = 1.05453547094965
$ws.Cells.Item(20, 9).Value = 1.040265794986932
$ws.Cells.Item(20, 10).Value = 1.043408503103173
$ws.Cells.Item(20, 11).Value = 1.044367024210032
$ws.Cells.Item(20, 12).Value = 1.048718533186343
$ws.Cells.Item(20, 13).Value = 1.05784550881811
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.036245353165363
$ws.Cells.Item(21, 4).Value = 1.040273365209623
$ws.Cells.Item(21, 5).Value = 1.044500739898518
$ws.Cells.Item(21, 6).Value = 1.053538652979852
$ws.Cells.Item(21, 9).Value = 1.040024656133433
$ws.Cells.Item(21, 10).Value = 1.042829382912474
$ws.Cells.Item(21, 11).Value = 1.043837671828035
$ws.Cells.Item(21, 12).Value = 1.04804952628671
$ws.Cells.Item(21, 13).Value = 1.057054719158973
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.035632896038937
$ws.Cells.Item(22, 4).Value = 1.039810228095669
$ws.Cells.Item(22, 5).Value = 1.04395061013062
$ws.Cells.Item(22, 6).Value = 1.052913197289821
$ws.Cells.Item(22, 9).Value = 1.039872162049297
$ws.Cells.Item(22, 10).Value = 1.04246527947795
$ws.Cells.Item(22, 11).Value = 1.04350468113699
$ws.Cells.Item(22, 12).Value = 1.047629261140086
$ws.Cells.Item(22, 13).Value = 1.05655811337287
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.035957502139895
$ws.Cells.Item(23, 4).Value = 1.040055690472105
$ws.Cells.Item(23, 5).Value = 1.044242141390608
$ws.Cells.Item(23, 6).Value = 1.053244640916546
$ws.Cells.Item(23, 9).Value = 1.039953085281196
$ws.Cells.Item(23, 10).Value = 1.042658296007871
$ws.Cells.Item(23, 11).Value = 1.043681221108257
$ws.Cells.Item(23, 12).Value = 1.047852016481386
$ws.Cells.Item(23, 13).Value = 1.056821316563448
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.037236967766223
$ws.Cells.Item(24, 4).Value = 1.041023270965578
$ws.Cells.Item(24, 5).Value = 1.045392145334214
$ws.Cells.Item(24, 6).Value = 1.054552204352469
$ws.Cells.Item(24, 9).Value = 1.040269822734394
$ws.Cells.Item(24, 10).Value = 1.043418212249163
$ws.Cells.Item(24, 11).Value = 1.044375895970075
$ws.Cells.Item(24, 12).Value = 1.048729755277535
$ws.Cells.Item(24, 13).Value = 1.057858776482439
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.03872481187598
$ws.Cells.Item(25, 4).Value = 1.042148553380381
$ws.Cells.Item(25, 5).Value = 1.046731302121332
$ws.Cells.Item(25, 6).Value = 1.056075084277759
$ws.Cells.Item(25, 9).Value = 1.040633521760836
$ws.Cells.Item(25, 10).Value = 1.0443000652751
$ws.Cells.Item(25, 11).Value = 1.045181263722694
$ws.Cells.Item(25, 12).Value = 1.049749869112519
$ws.Cells.Item(25, 13).Value = 1.059065228961354
